$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Name/Email values between row 3 (was guru) and row 4 (was alekhya)
$ws.Range("B3").Value = "alekhya"
$ws.Range("C3").Value = "alekhyakanjarla@gmail.com"
$ws.Range("B4").Value = "guru"
$ws.Range("C4").Value = "guru.sai.shreesh@gmail.com"

# The hyperlink that lives on C3 still points at guru.sai.shreesh@gmail.com;
# give it an explicit display text matching that target (foreach enumeration
# binds correctly to the live hyperlink object, unlike indexed .Item access)
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -eq "mailto:guru.sai.shreesh@gmail.com") {
        $hl.TextToDisplay = "guru.sai.shreesh@gmail.com"
    }
}

# Update the active selection to A6
$ws.Range("A6").Select()
